# Weekly update: insert a new week's price block for "Lechuga" (Femacal de La Calera)
# above the existing data, shifting the old rows down by 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows at the top of the data block (old row 906 onward shifts down to 911).
$ws.Rows("906:910").Insert()

# Common values for the whole new block (same as the rest of this sub-sheet).
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$fecha     = 44448
$codreg    = 5
$catId     = 100112033
$categoria = "Lechuga"
$origen    = "Provincia de Quillota"
$clasif    = "Hortaliza"

# Per-variety values: row, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Unidad, Precio$/Kg, KgOUnidades
$data = @(
    @(906, "Conconina(o)",    "Primera", 160, 5500, 6000, 5750, "`$/caja 10 unidades", 575, 10),
    @(907, "Escarola",        "Primera", 175, 7500, 8000, 7743, "`$/caja 15 unidades", 516, 15),
    @(908, "Francesa morada", "Primera", 140, 6000, 6300, 6171, "`$/caja 18 unidades", 343, 18),
    @(909, "Marina",          "Primera", 130, 6000, 6500, 6231, "`$/caja 18 unidades", 346, 18),
    @(910, "Milanesa",        "Primera", 130, 5500, 5800, 5662, "`$/caja 20 unidades", 283, 20)
)

foreach ($row in $data) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $catId
    $ws.Cells.Item($r, 7).Value  = $categoria
    $ws.Cells.Item($r, 8).Value  = $row[1]
    $ws.Cells.Item($r, 9).Value  = $row[2]
    $ws.Cells.Item($r, 10).Value = $row[3]
    $ws.Cells.Item($r, 11).Value = $row[4]
    $ws.Cells.Item($r, 12).Value = $row[5]
    $ws.Cells.Item($r, 13).Value = $row[6]
    $ws.Cells.Item($r, 14).Value = $row[7]
    $ws.Cells.Item($r, 15).Value = $origen
    $ws.Cells.Item($r, 16).Value = $row[8]
    $ws.Cells.Item($r, 17).Value = $row[9]
    $ws.Cells.Item($r, 18).Value = $clasif
}
